# Added test data for Switzerland market.
#
# - Select all cells on the "Germany" sheet (mirrors the author reviewing
#   it before switching away - it ends up with a full-sheet selection and
#   loses the "active tab" flag once focus moves elsewhere).
# - Duplicate the "Czech" sheet (same layout/styles) to a new sheet placed
#   right after it, rename it "Swiss", and update the two cells that are
#   specific to the market: the market name and the NGC/Ticket reference.
# - Finish with the selection on B4 of the new sheet, which becomes the
#   active tab.

$wb = $excel.ActiveWorkbook

$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Cells.Select()

$czech = $wb.Worksheets.Item("Czech")
$czech.Copy([System.Reflection.Missing]::Value, $czech)

# The copy is inserted right after "Czech" and becomes the last sheet.
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2653"

$swiss.Range("B4").Select()
